$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy'
$ws.Range('G4').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud'
$ws.Range('G5').Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Nada Gouda'
$ws.Range('G6').Value = 'Dr. Kerelos Zareef, Dr. Nada Mohammad'
$ws.Range('G10').Value = 'Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Marina Youhanna, Dr. Maryam Ahmad'
$ws.Range('G12').Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range('G18').Value = 'Dr. Yasmin, Dr. Shorok Mohammad, Dr. Remon, Dr. Aya Hanafy'
$ws.Range('G19').Value = 'Dr. Nardine, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Monica, Dr. Salma Hassan, Dr. Maryam Ashraf, Dr. Remon'
$ws.Range('G20').Value = 'Dr. Nardine, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Aya Hanafy, Dr. Youstina Magdy, Dr. Marina Sorial, Dr. Remon'
$ws.Range('G21').Value = 'Dr. Monica, Dr. Yasmin, Dr. Yassmen Ahmad, Dr. Shorok Mohammad'
$ws.Range('G24').Value = 'Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy'
$ws.Range('G25').Value = 'Dr. Alshimaa Atef, Administrator, Dr. Manar Montaser, Dr. Gehan Adel'
$ws.Range('G26').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Hend Mahmoud'
$ws.Range('G27').Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Nada Gouda'
$ws.Range('G28').Value = 'Dr. Kerelos Zareef, Dr. Nada Mohammad'
$ws.Range('G32').Value = 'Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Marina Youhanna, Dr. Maryam Ahmad'
$ws.Range('G34').Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range('G40').Value = 'Dr. Yasmin, Dr. Shorok Mohammad, Dr. Remon, Dr. Aya Hanafy'
$ws.Range('G41').Value = 'Dr. Nardine, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Monica, Dr. Salma Hassan, Dr. Maryam Ashraf, Dr. Remon'
$ws.Range('G42').Value = 'Dr. Nardine, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Aya Hanafy, Dr. Youstina Magdy, Dr. Marina Sorial, Dr. Remon'
$ws.Range('G43').Value = 'Dr. Monica, Dr. Yasmin, Dr. Yassmen Ahmad, Dr. Shorok Mohammad'
$ws.Range('G48').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub'
$ws.Range('G54').Value = 'Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Mai Mustafa'
$ws.Range('G62').Value = 'Dr. Shorok Mohammad, Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida'
$ws.Range('G63').Value = 'Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah'
$ws.Range('G64').Value = 'Dr. Wafaa Ebida, Dr. Youstina Magdy'
$ws.Range('G65').Value = 'Dr. Ola Abd Al-Fattah, Dr. Nardine, Dr. Aya Hanafy, Dr. Eman Samir Gabry, Dr. Shorok Mohammad, Dr. Salma Hassan, Dr. Remon'
$ws.Range('G66').Value = 'Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Monica, Dr. Marina Sorial, Dr. Maryam Ashraf'
$ws.Range('G70').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub'
$ws.Range('G76').Value = 'Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Mai Mustafa'
$ws.Range('G81').Value = 'Dr. Enas Omran, Dr. Walaa Ghanima'
$ws.Range('G84').Value = 'Dr. Shorok Mohammad, Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida'
$ws.Range('G85').Value = 'Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Ola Abd Al-Fattah'
$ws.Range('G86').Value = 'Dr. Wafaa Ebida, Dr. Youstina Magdy'
$ws.Range('G87').Value = 'Dr. Ola Abd Al-Fattah, Dr. Nardine, Dr. Aya Hanafy, Dr. Eman Samir Gabry, Dr. Shorok Mohammad, Dr. Salma Hassan, Dr. Remon'
$ws.Range('G88').Value = 'Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Monica, Dr. Marina Sorial, Dr. Maryam Ashraf'
$ws.Range('G90').Value = 'Dr. Mohammad El-Tanany, Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki'
$ws.Range('G92').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub'
$ws.Range('G93').Value = 'Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Amera Ahmad Saad'
$ws.Range('G96').Value = 'Dr. Mariam Nour El-Din, Dr. Sara Nabil, Dr. Nourhan Mohammad, Dr. Amal Awwad'
$ws.Range('G98').Value = 'Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Mai Mustafa'
$ws.Range('G106').Value = 'Dr. Nardine, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Monica, Dr. Remon'
$ws.Range('G107').Value = 'Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf'
$ws.Range('G108').Value = 'Dr. Nardine, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Aya Hanafy, Dr. Youstina Magdy, Dr. Marina Sorial, Dr. Remon'
$ws.Range('G111').Value = 'Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Naema Gomaa, Dr. Monica, Dr. Eman Samir Gabry, Dr. Marina Atef'
$ws.Range('G112').Value = 'Dr. Mohammad El-Tanany, Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki'
$ws.Range('G114').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub'
$ws.Range('G115').Value = 'Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Menna tu''Alllah Mohammad, Dr. Amera Ahmad Saad'
$ws.Range('G118').Value = 'Dr. Mariam Nour El-Din, Dr. Sara Nabil, Dr. Nourhan Mohammad, Dr. Amal Awwad'
$ws.Range('G120').Value = 'Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Merna Said, Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Maryam Ahmad, Dr. Mai Mustafa'
$ws.Range('G128').Value = 'Dr. Nardine, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Monica, Dr. Remon'
$ws.Range('G129').Value = 'Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf'
$ws.Range('G130').Value = 'Dr. Nardine, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Aya Hanafy, Dr. Youstina Magdy, Dr. Marina Sorial, Dr. Remon'
$ws.Range('G133').Value = 'Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Naema Gomaa, Dr. Monica, Dr. Eman Samir Gabry, Dr. Marina Atef'
$ws.Range('G134').Value = 'Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Veronia Rafat'
$ws.Range('G137').Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Nada Gouda'
$ws.Range('G142').Value = 'Dr. Marwa Mustafa, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Basma Hamed, Dr. Yasmeena Fattoh'
$ws.Range('G150').Value = 'Dr. Nardine, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Monica, Dr. Salma Hassan, Dr. Maryam Ashraf, Dr. Remon'
$ws.Range('G151').Value = 'Dr. Monica, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Marina Atef'
$ws.Range('G152').Value = 'Dr. Wafaa Ebida, Dr. Marina Atef'
$ws.Range('G153').Value = 'Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Monica, Dr. Marina Sorial, Dr. Maryam Ashraf'
$ws.Range('G154').Value = 'Dr. Remon, Dr. Naema Gomaa, Dr. Wafaa Ebida, Dr. Salma Hassan'
$ws.Range('G155').Value = 'Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Naema Gomaa, Dr. Monica, Dr. Eman Samir Gabry, Dr. Marina Atef'
$ws.Range('G156').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef, Dr. Majorelle Magdy, Dr. Manar Montaser'
$ws.Range('G159').Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Fatma Elhady, Dr. Nada Gouda'
$ws.Range('G164').Value = 'Dr. Marwa Mustafa, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Basma Hamed, Dr. Yasmeena Fattoh'
$ws.Range('G165').Value = 'Dr. Sarah Mahdy, Dr. Nouran Mahmoud'
$ws.Range('G172').Value = 'Dr. Nardine, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Monica, Dr. Salma Hassan, Dr. Maryam Ashraf, Dr. Remon'
$ws.Range('G173').Value = 'Dr. Monica, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Marina Atef'
$ws.Range('G174').Value = 'Dr. Wafaa Ebida, Dr. Marina Atef'
$ws.Range('G175').Value = 'Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Monica, Dr. Marina Sorial, Dr. Maryam Ashraf'
$ws.Range('G176').Value = 'Dr. Remon, Dr. Naema Gomaa, Dr. Wafaa Ebida, Dr. Salma Hassan'
$ws.Range('G177').Value = 'Dr. Ola Abd Al-Fattah, Dr. Yasmin, Dr. Naema Gomaa, Dr. Monica, Dr. Eman Samir Gabry, Dr. Marina Atef'
